# Update the "Estado de Cuenta" worker data table.
# Rows 16-17 swap their worker identity (Duvan Caicedo moves up, Orlando
# Ordoñez moves down and changes period), row 18 becomes Orlando Ordoñez's
# second record for period 2110. Row 19 (Jose Carlos Hernandez) is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16: Duvan Caicedo Gaviria, period 2110
$ws.Range("C16").Value = "1238338773"
$ws.Range("D16").Value = "DUVAN CAICEDO GAVIRIA"
$ws.Range("E16").Value = "2110"
$ws.Range("F16").Value = 23016
$ws.Range("G16").Value = 908526

# Row 17: Orlando Ordoñez Gonzales, period 2111
$ws.Range("C17").Value = "91101519"
$ws.Range("D17").Value = "ORLANDO ORDOÑEZ GONZALES"
$ws.Range("E17").Value = "2111"
$ws.Range("F17").Value = 38422
$ws.Range("G17").Value = 960528

# Row 18: Orlando Ordoñez Gonzales, period 2110
$ws.Range("C18").Value = "91101519"
$ws.Range("D18").Value = "ORLANDO ORDOÑEZ GONZALES"
$ws.Range("E18").Value = "2110"
$ws.Range("F18").Value = 24334
$ws.Range("G18").Value = 960528
